$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 24.832978138124986
$ws.Range("C2").Value = 23.28494295750005
$ws.Range("D2").Value = 31.556935828125006
$ws.Range("E2").Value = 34.665075633125014

$ws.Range("B3").Value = 20.167593688124953
$ws.Range("C3").Value = 38.205168322500015
$ws.Range("D3").Value = 35.443980539999927
$ws.Range("E3").Value = 31.535374552500002

$ws.Range("B1:E3").Select()
